$d = $word.ActiveDocument

$replacements = @(
    @{old='651×9=5859'; new='971×2=1942'},
    @{old='839×2=1678'; new='569×9=5121'},
    @{old='232×3=696';  new='732×4=2928'},
    @{old='823×8=6584'; new='851×9=7659'},
    @{old='852×2=1704'; new='830×8=6640'},
    @{old='521×2=1042'; new='936×8=7488'},
    @{old='285×5=1425'; new='345×6=2070'},
    @{old='957×5=4785'; new='195×8=1560'},
    @{old='934×2=1868'; new='569×6=3414'},
    @{old='502×6=3012'; new='390×9=3510'},
    @{old='244×2=488';  new='918×2=1836'},
    @{old='368×3=1104'; new='686×6=4116'},
    @{old='221×3=663';  new='383×4=1532'},
    @{old='461×4=1844'; new='562×6=3372'},
    @{old='443×9=3987'; new='247×3=741'},
    @{old='678×4=2712'; new='334×3=1002'},
    @{old='618×7=4326'; new='459×8=3672'},
    @{old='258×7=1806'; new='898×4=3592'},
    @{old='173×9=1557'; new='773×3=2319'},
    @{old='157×7=1099'; new='246×6=1476'},
    @{old='465×9=4185'; new='477×8=3816'},
    @{old='967×8=7736'; new='830×4=3320'},
    @{old='690×8=5520'; new='909×6=5454'},
    @{old='947×5=4735'; new='267×4=1068'},
    @{old='182×8=1456'; new='730×2=1460'}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
